$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Delete the unused "Registration" column (column I)
$ws.Columns("I").Delete()

# Delete the now-orphaned note row (was row 27) left over after the column shift
$ws.Rows("27").Delete()
